# Update crypto price/volume data as scraped on the latest GitHub Actions run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("Price") holds numeric-looking values that must stay as plain text
# (e.g. "587.32", "1.00", thousand-separated "69.479.40"), so force text format
# before assigning them - otherwise Excel auto-converts them to real numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.479.40"
$ws.Range("E2").Value = "  +2.35%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.391.68"
$ws.Range("E3").Value = "  +1.87%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "587.32"
$ws.Range("E5").Value = "  +0.83%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "180.15"
$ws.Range("E6").Value = "  +1.74%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.596"
$ws.Range("E8").Value = "  +0.88%  "

$ws.Range("E9").Value = "  +5.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.592"
$ws.Range("E10").Value = "  +1.53%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "48.54"
$ws.Range("E11").Value = "  +2.19%  "

$ws.Range("E12").Value = "  +3.23%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "680.33"
$ws.Range("E13").Value = "  -3.14%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "8.63"
$ws.Range("E14").Value = "  +2.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.935.14"
$ws.Range("E15").Value = "  +1.71%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.489.89"
$ws.Range("E16").Value = "  +2.31%  "

$ws.Range("E17").Value = "  +1.73%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.389.16"
$ws.Range("E18").Value = "  +1.69%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "17.64"
$ws.Range("E19").Value = "  +0.80%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.29"
$ws.Range("E20").Value = "  +1.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.904"
$ws.Range("E21").Value = "  +0.85%  "

$ws.Range("E22").Value = "  +0.56%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "17.21"
$ws.Range("E23").Value = "  +0.68%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "103.37"
$ws.Range("E24").Value = "  +3.75%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "3.94"
$ws.Range("E25").Value = "  +0.56%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.73"
$ws.Range("E26").Value = "  +1.16%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.70"
$ws.Range("E27").Value = "  +1.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "34.10"
$ws.Range("E28").Value = "  +3.01%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.71"
$ws.Range("E29").Value = "  +1.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.08"
$ws.Range("E30").Value = "  -0.10%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "11.14"
$ws.Range("E31").Value = "  +1.03%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "557.57"
$ws.Range("E32").Value = "  -1.88%  "

$ws.Range("E33").Value = "  +5.05%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.107"
$ws.Range("E34").Value = "  +0.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "58.33"
$ws.Range("E35").Value = "  +1.70%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.682.71"
$ws.Range("E37").Value = "  -0.06%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "35.47"
$ws.Range("E38").Value = "  +2.22%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.139"
$ws.Range("E39").Value = "  +4.05%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.28"
$ws.Range("E40").Value = "  +3.51%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.72"
$ws.Range("E41").Value = "  +2.08%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0₃0698"
$ws.Range("E42").Value = "  +2.99%  "

$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0424"
$ws.Range("E44").Value = "  +3.99%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.29"
$ws.Range("E45").Value = "  -1.32%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.68"
$ws.Range("E46").Value = "  +0.10%  "

$ws.Range("E47").Value = "  +0.81%  "

$ws.Range("E48").Value = "  +5.54%  "

$ws.Range("E49").Value = "  -0.05%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "133.20"
$ws.Range("E50").Value = "  +1.53%  "

$ws.Range("E51").Value = "  +3.78%  "
